$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1: "- Flags. ..." paragraph - split the sentence and append new text
# about flags being mutually exclusive.
# ---------------------------------------------------------------------------
$flagsFind = $d.Content
$flagsFind.Find.Execute("- Flags. These fields serve to indicate") | Out-Null
$flagsPara = $flagsFind.Paragraphs(1).Range

$flagsXml = @"
<w:p $wns>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">- Flags. These fields serve to indicate in the code to which group the material belongs in terms of functionality within the heating system. </w:t>
  </w:r>
  <w:r>
    <w:t>A material can either be invariant, thermoelectric, a phase-change material, magnetocaloric, electrocaloric, elastocaloric or barocaloric. At the moment, these flags are mutually exclusive. Should you have a material, that is multicaloric, make two or more materials with same IDs, with different appended letters. However, if you wish to use multicaloric effect in a simulation simultaneously, open an issue and we will implement that.</w:t>
  </w:r>
</w:p>
"@
$flagsPara.InsertXML($flagsXml)

# ---------------------------------------------------------------------------
# Change 2: restructure the "- The info.json file ..." paragraph through the
# hysteresis paragraph: merge some runs, split the final paragraph into four
# separate paragraphs, and move the lastRenderedPageBreak + reword a couple
# of phrases.
# ---------------------------------------------------------------------------
$startFind = $d.Content
$startFind.Find.Execute("- The info.json file") | Out-Null
$startPara = $startFind.Paragraphs(1).Range

$endFind = $d.Content
$endFind.Find.Execute("- There are also other options where the mentioned three files") | Out-Null
$endPara = $endFind.Paragraphs(1).Range

$blockRange = $d.Range($startPara.Start, $endPara.End)

$blockXml = @"
<w:p $wns w:rsidR="00B33016" w:rsidRDefault="00B33016" w:rsidP="00B33016">
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t>- The info.</w:t>
  </w:r>
  <w:r w:rsidR="004F3109">
    <w:t>json</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> file </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">is in JSON format and must </w:t>
  </w:r>
  <w:r>
    <w:t>contain</w:t>
  </w:r>
  <w:r>
    <w:t>: RT properties, ranges, fields, ID, short name, long name, and melting point. The RT properties must contain the values of density, specific heat, conductivity and emissivity at room temperature. The ranges must contain the temperature ranges in which each material property is defined. Temperature ranges must be recorded for density, specific heat capacity, thermal conductivity, adiabatic temperature change (if the material is caloric), emissivity, and other relevant properties (e.g. seebeck coefficient for thermoelectric materials, etc.). The fields contains strengths of the external fields where properties for caloric materials are defined. When it comes to magnetic fields, the values are in T, when it comes to electric fields, they are in MV/m, and for pressure and stress, they are in bars. (See any material for example).</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t>- Files rho.txt, cp.txt and k.txt (either one value at room temperature, or one column of 20000 values from 0 to 2000 K in steps of 0.1 K); here there are only values without temperatures.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>- One or more of the above three files can be replaced by several files for different external fields (magnetic</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>written in T, electric in MVm (which means MV/m)), e.g. cp_0.0T.txt, cp_1.0T.txt, etc. This happens e.g. in caloric materials.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t>- There are also other options where the mentioned three files are replaced by some hysteresis, e.g. cp_heating.txt and cp_cooling.txt, i.e. when the material</w:t>
  </w:r>
  <w:r>
    <w:t>’s</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> cpThysteresis flag</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> is</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> true. But there can also be hysteresis with different fields, then we get cp_0.0T_cooling.txt, cp_0.0T_heating.txt, cp_1.0T_cooling.txt, etc.</w:t>
  </w:r>
</w:p>
"@
$blockRange.InsertXML($blockXml)
